$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 46065
$ws.Range("C3").Value = 46065
$ws.Range("C4").Value = 46065
$ws.Range("C5").Value = 46065

$ws.Range("A6").Value = "A 19003-2025"
$ws.Range("B6").Value = 45764
$ws.Range("C6").Value = 46065
$ws.Range("G6").Value = 5.4

$ws.Range("A7").Value = "A 1468-2022"
$ws.Range("B7").Value = 44573
$ws.Range("C7").Value = 46065
$ws.Range("G7").Value = 1.8

$ws.Range("C8").Value = 46065

$ws.Range("A9").Value = "A 18968-2025"
$ws.Range("B9").Value = 45764.53686342593
$ws.Range("C9").Value = 46065
$ws.Range("G9").Value = 8.199999999999999

$ws.Range("A10").Value = "A 65018-2023"
$ws.Range("B10").Value = 45287
$ws.Range("C10").Value = 46065
$ws.Range("G10").Value = 1.1

$ws.Range("A11").Value = "A 53361-2024"
$ws.Range("B11").Value = 45614
$ws.Range("C11").Value = 46065
$ws.Range("G11").Value = 2.5

$ws.Range("A12").Value = "A 53361-2024"
$ws.Range("B12").Value = 45614
$ws.Range("C12").Value = 46065
$ws.Range("G12").Value = 0.4

$ws.Range("A13").Value = "A 53750-2025"
$ws.Range("B13").Value = 45960.65806712963
$ws.Range("C13").Value = 46065
$ws.Range("G13").Value = 0.9
